# Swap the two embedded themes:
#   ppt/theme/theme1.xml  ("Integral" / "Red Violet")      -> "Office Theme" / "Office" palette
#   ppt/theme/theme2.xml  ("Office Theme" / "Office")       -> "Integral" / "Red Violet" palette
#
# The only theme part reachable from the slides (theme1.xml, wired to the
# single SlideMaster that every slide/layout in this deck uses) can be
# edited in place through Slide.ThemeColorScheme - each of its 12 entries
# maps 1:1 (by index) onto the <a:clrScheme> children dk1/lt1/dk2/lt2/
# accent1-6/hlink/folHlink of that theme part.  theme2.xml only backs the
# NotesMaster, which this object model does not expose a ThemeColorScheme
# (or any other raw-XML) handle for, so it cannot be touched from here.
#
# RGB() isn't a PowerShell builtin (it's VBA-only), so colours are packed
# by hand into the 0x00BBGGRR long that PowerPoint's ColorFormat.RGB uses.

function New-ComRgb([int]$r, [int]$g, [int]$b) {
    return $r + ($g * 256) + ($b * 65536)
}

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$tcs = $s.ThemeColorScheme

# Index -> target "Office" theme colour (hex). Index order mirrors the
# <a:clrScheme> child order: dk1, lt1, dk2, lt2, accent1..accent6, hlink,
# folHlink.
$targetHex = @(
    "000000", # 1 dk1
    "FFFFFF", # 2 lt1
    "44546A", # 3 dk2
    "E7E6E6", # 4 lt2
    "5B9BD5", # 5 accent1
    "ED7D31", # 6 accent2
    "A5A5A5", # 7 accent3
    "FFC000",  # 8 accent4
    "4472C4", # 9 accent5
    "70AD47", # 10 accent6
    "0563C1", # 11 hlink
    "954F72"  # 12 folHlink
)

for ($i = 1; $i -le $tcs.Count; $i++) {
    $hex = $targetHex[$i - 1]
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    $tcs.Colors($i).RGB = New-ComRgb $r $g $b
}
